$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.323799641721889
$ws.Range("D2").Value = 0.7555442988502048
$ws.Range("E2").Value = -0.2498931886842296
$ws.Range("F2").Value = -0.006045411213057164
$ws.Range("H2").Value = 0.323799641721889
$ws.Range("I2").Value = 0.7555442988502048
$ws.Range("J2").Value = -0.2498931886842296
$ws.Range("K2").Value = -0.006045411213057164
$ws.Range("M2").Value = -0.01917933785190364
$ws.Range("N2").Value = -0.5864791365772161
$ws.Range("O2").Value = 0.2019530246563825
$ws.Range("P2").Value = -0.2545595402151047
$ws.Range("B3").Value = 0.323799641721889
$ws.Range("D3").Value = 0.3962439059502366
$ws.Range("E3").Value = -0.1800546013261557
$ws.Range("F3").Value = -0.3362354095419539
$ws.Range("G3").Value = 0.3237996417218891
$ws.Range("H3").Value = 1.0
$ws.Range("I3").Value = 0.3962439059502366
$ws.Range("J3").Value = -0.1800546013261557
$ws.Range("K3").Value = -0.3362354095419539
$ws.Range("L3").Value = 0.3237996417218891
$ws.Range("M3").Value = 0.4130445827457469
$ws.Range("N3").Value = -0.4611548869416665
$ws.Range("O3").Value = 0.2425999198641816
$ws.Range("P3").Value = 0.07364185281972291
$ws.Range("B4").Value = 0.7555442988502048
$ws.Range("C4").Value = 0.3962439059502366
$ws.Range("E4").Value = -0.1426760180749146
$ws.Range("F4").Value = -0.03264241957040208
$ws.Range("G4").Value = 0.7555442988502049
$ws.Range("H4").Value = 0.3962439059502366
$ws.Range("I4").Value = 1.0
$ws.Range("J4").Value = -0.1426760180749146
$ws.Range("K4").Value = -0.03264241957040208
$ws.Range("L4").Value = 0.7555442988502049
$ws.Range("M4").Value = 0.09237240891518564
$ws.Range("N4").Value = -0.5471734752080594
$ws.Range("O4").Value = 0.2016028559016793
$ws.Range("P4").Value = -0.1892737362199487
$ws.Range("B5").Value = -0.2498931886842296
$ws.Range("C5").Value = -0.1800546013261557
$ws.Range("D5").Value = -0.1426760180749146
$ws.Range("F5").Value = 0.02530208406707743
$ws.Range("G5").Value = -0.2498931886842295
$ws.Range("H5").Value = -0.1800546013261557
$ws.Range("I5").Value = -0.1426760180749146
$ws.Range("J5").Value = 1.0
$ws.Range("K5").Value = 0.02530208406707743
$ws.Range("L5").Value = -0.2498931886842295
$ws.Range("M5").Value = 0.06722331682273598
$ws.Range("N5").Value = 0.1272006909638258
$ws.Range("O5").Value = -0.04649015646051528
$ws.Range("P5").Value = 0.03693963816442876
$ws.Range("B6").Value = -0.006045411213057164
$ws.Range("C6").Value = -0.3362354095419539
$ws.Range("D6").Value = -0.03264241957040208
$ws.Range("E6").Value = 0.02530208406707743
$ws.Range("G6").Value = -0.006045411213057154
$ws.Range("H6").Value = -0.336235409541954
$ws.Range("I6").Value = -0.03264241957040209
$ws.Range("J6").Value = 0.02530208406707743
$ws.Range("K6").Value = 1.0
$ws.Range("L6").Value = -0.006045411213057154
$ws.Range("M6").Value = -0.545766696202424
$ws.Range("N6").Value = 0.4941076255294825
$ws.Range("O6").Value = -0.4774024123023282
$ws.Range("P6").Value = -0.6169364305255066
$ws.Range("C7").Value = 0.3237996417218891
$ws.Range("D7").Value = 0.7555442988502049
$ws.Range("E7").Value = -0.2498931886842295
$ws.Range("F7").Value = -0.006045411213057154
$ws.Range("H7").Value = 0.323799641721889
$ws.Range("I7").Value = 0.7555442988502048
$ws.Range("J7").Value = -0.2498931886842296
$ws.Range("K7").Value = -0.006045411213057164
$ws.Range("M7").Value = -0.01917933785190364
$ws.Range("N7").Value = -0.5864791365772161
$ws.Range("O7").Value = 0.2019530246563825
$ws.Range("P7").Value = -0.2545595402151047
$ws.Range("B8").Value = 0.323799641721889
$ws.Range("C8").Value = 1.0
$ws.Range("D8").Value = 0.3962439059502366
$ws.Range("E8").Value = -0.1800546013261557
$ws.Range("F8").Value = -0.336235409541954
$ws.Range("G8").Value = 0.323799641721889
$ws.Range("I8").Value = 0.3962439059502366
$ws.Range("J8").Value = -0.1800546013261557
$ws.Range("K8").Value = -0.3362354095419539
$ws.Range("L8").Value = 0.3237996417218891
$ws.Range("M8").Value = 0.4130445827457469
$ws.Range("N8").Value = -0.4611548869416665
$ws.Range("O8").Value = 0.2425999198641816
$ws.Range("P8").Value = 0.07364185281972291
$ws.Range("B9").Value = 0.7555442988502048
$ws.Range("C9").Value = 0.3962439059502366
$ws.Range("D9").Value = 1.0
$ws.Range("E9").Value = -0.1426760180749146
$ws.Range("F9").Value = -0.03264241957040209
$ws.Range("G9").Value = 0.7555442988502048
$ws.Range("H9").Value = 0.3962439059502366
$ws.Range("J9").Value = -0.1426760180749146
$ws.Range("K9").Value = -0.03264241957040208
$ws.Range("L9").Value = 0.7555442988502049
$ws.Range("M9").Value = 0.09237240891518564
$ws.Range("N9").Value = -0.5471734752080594
$ws.Range("O9").Value = 0.2016028559016793
$ws.Range("P9").Value = -0.1892737362199487
$ws.Range("B10").Value = -0.2498931886842296
$ws.Range("C10").Value = -0.1800546013261557
$ws.Range("D10").Value = -0.1426760180749146
$ws.Range("E10").Value = 1.0
$ws.Range("F10").Value = 0.02530208406707743
$ws.Range("G10").Value = -0.2498931886842296
$ws.Range("H10").Value = -0.1800546013261557
$ws.Range("I10").Value = -0.1426760180749146
$ws.Range("K10").Value = 0.02530208406707743
$ws.Range("L10").Value = -0.2498931886842295
$ws.Range("M10").Value = 0.06722331682273598
$ws.Range("N10").Value = 0.1272006909638258
$ws.Range("O10").Value = -0.04649015646051528
$ws.Range("P10").Value = 0.03693963816442876
$ws.Range("B11").Value = -0.006045411213057164
$ws.Range("C11").Value = -0.3362354095419539
$ws.Range("D11").Value = -0.03264241957040208
$ws.Range("E11").Value = 0.02530208406707743
$ws.Range("F11").Value = 1.0
$ws.Range("G11").Value = -0.006045411213057164
$ws.Range("H11").Value = -0.3362354095419539
$ws.Range("I11").Value = -0.03264241957040208
$ws.Range("J11").Value = 0.02530208406707743
$ws.Range("L11").Value = -0.006045411213057154
$ws.Range("M11").Value = -0.545766696202424
$ws.Range("N11").Value = 0.4941076255294825
$ws.Range("O11").Value = -0.4774024123023282
$ws.Range("P11").Value = -0.6169364305255066
$ws.Range("C12").Value = 0.3237996417218891
$ws.Range("D12").Value = 0.7555442988502049
$ws.Range("E12").Value = -0.2498931886842295
$ws.Range("F12").Value = -0.006045411213057154
$ws.Range("H12").Value = 0.3237996417218891
$ws.Range("I12").Value = 0.7555442988502049
$ws.Range("J12").Value = -0.2498931886842295
$ws.Range("K12").Value = -0.006045411213057154
$ws.Range("M12").Value = -0.01917933785190364
$ws.Range("N12").Value = -0.5864791365772161
$ws.Range("O12").Value = 0.2019530246563825
$ws.Range("P12").Value = -0.2545595402151047
$ws.Range("B13").Value = -0.01917933785190364
$ws.Range("C13").Value = 0.4130445827457469
$ws.Range("D13").Value = 0.09237240891518564
$ws.Range("E13").Value = 0.06722331682273598
$ws.Range("F13").Value = -0.545766696202424
$ws.Range("G13").Value = -0.01917933785190364
$ws.Range("H13").Value = 0.4130445827457469
$ws.Range("I13").Value = 0.09237240891518564
$ws.Range("J13").Value = 0.06722331682273598
$ws.Range("K13").Value = -0.545766696202424
$ws.Range("L13").Value = -0.01917933785190364
$ws.Range("N13").Value = 0.04820155874753494
$ws.Range("O13").Value = 0.8780838429650567
$ws.Range("P13").Value = -0.08138554629115878
$ws.Range("B14").Value = -0.5864791365772161
$ws.Range("C14").Value = -0.4611548869416665
$ws.Range("D14").Value = -0.5471734752080594
$ws.Range("E14").Value = 0.1272006909638258
$ws.Range("F14").Value = 0.4941076255294825
$ws.Range("G14").Value = -0.5864791365772161
$ws.Range("H14").Value = -0.4611548869416665
$ws.Range("I14").Value = -0.5471734752080594
$ws.Range("J14").Value = 0.1272006909638258
$ws.Range("K14").Value = 0.4941076255294825
$ws.Range("L14").Value = -0.5864791365772161
$ws.Range("M14").Value = 0.04820155874753494
$ws.Range("O14").Value = 0.007524513946420188
$ws.Range("P14").Value = -0.536789932940771
$ws.Range("B15").Value = 0.2019530246563825
$ws.Range("C15").Value = 0.2425999198641816
$ws.Range("D15").Value = 0.2016028559016793
$ws.Range("E15").Value = -0.04649015646051528
$ws.Range("F15").Value = -0.4774024123023282
$ws.Range("G15").Value = 0.2019530246563825
$ws.Range("H15").Value = 0.2425999198641816
$ws.Range("I15").Value = 0.2016028559016793
$ws.Range("J15").Value = -0.04649015646051528
$ws.Range("K15").Value = -0.4774024123023282
$ws.Range("L15").Value = 0.2019530246563825
$ws.Range("M15").Value = 0.8780838429650567
$ws.Range("N15").Value = 0.007524513946420188
$ws.Range("P15").Value = -0.2182122118000741
$ws.Range("B16").Value = -0.2545595402151047
$ws.Range("C16").Value = 0.07364185281972291
$ws.Range("D16").Value = -0.1892737362199487
$ws.Range("E16").Value = 0.03693963816442876
$ws.Range("F16").Value = -0.6169364305255066
$ws.Range("G16").Value = -0.2545595402151047
$ws.Range("H16").Value = 0.07364185281972291
$ws.Range("I16").Value = -0.1892737362199487
$ws.Range("J16").Value = 0.03693963816442876
$ws.Range("K16").Value = -0.6169364305255066
$ws.Range("L16").Value = -0.2545595402151047
$ws.Range("M16").Value = -0.08138554629115878
$ws.Range("N16").Value = -0.536789932940771
$ws.Range("O16").Value = -0.2182122118000741
